$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B2 currently holds the numeric value 45; replace it with the text date
# string "2000-07-02" (stored as text, formatted with the Text number
# format so Excel does not reinterpret it as a date/number).
$cell = $ws.Range("B2")
$cell.NumberFormat = "@"
$cell.Value = "2000-07-02"

# Column B is resized (best-fit) to accommodate the new date string.
$ws.Columns.Item(2).ColumnWidth = 10

# Move the active selection to I3, matching the post-edit selection state.
$ws.Range("I3").Select() | Out-Null
